$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, centered, bordered) onto the
# new I1/J1 header cells so they pick up the existing style index instead
# of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns (rows 2-11)
$values = @(7, 9, 9, 9, 9, 9, 9, 8, 9, 7)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
